$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.136.46"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "3.199.06"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.03%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").Value = "3.747.63"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.138"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("E14").Value = "  +2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "60.163.77"
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("D17").Value = "3.186.22"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.85%  "
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "0.0₃0903"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("E34").Value = "  +4.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").Value = "2.776.89"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0713"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.729"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0287"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.98%  "
$ws.Range("B45").Value = "RenzoRestakedETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D45").Value = "3.241.17"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.798"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("E51").Value = "  +0.03%  "

Write-Output "Applied all changes"
